$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume number, week-covering dates) ---
$ws.Range("A8").Value = "Volume 31   Number  31"
$ws.Range("C9").Value = "Report Covering the Week  7/29/2024  Through  8/4/2024"

# --- Cells changing from numeric to text ("0" / "***.*") ---
$ws.Range("C22").Copy($ws.Range("G15"))
$ws.Range("E22").Copy($ws.Range("H15"))
$ws.Range("C22").Copy($ws.Range("D16"))
$ws.Range("E22").Copy($ws.Range("E16"))
$ws.Range("C22").Copy($ws.Range("D23"))
$ws.Range("E22").Copy($ws.Range("E23"))
$ws.Range("C22").Copy($ws.Range("G27"))
$ws.Range("E22").Copy($ws.Range("H27"))

# --- Cells changing from text ("0"/"***.*") to numeric ---
$ws.Range("C15").Copy($ws.Range("C20"))
$ws.Range("C20").Value = 1
$ws.Range("C15").Copy($ws.Range("D28"))
$ws.Range("D28").Value = 1
$ws.Range("K16").Copy($ws.Range("E28"))
$ws.Range("E28").Value = 0
$ws.Range("C15").Copy($ws.Range("D29"))
$ws.Range("D29").Value = 1
$ws.Range("K16").Copy($ws.Range("E29"))
$ws.Range("E29").Value = -100
$ws.Range("C15").Copy($ws.Range("G29"))
$ws.Range("G29").Value = 1
$ws.Range("K16").Copy($ws.Range("H29"))
$ws.Range("H29").Value = -100
$ws.Range("C15").Copy($ws.Range("D30"))
$ws.Range("D30").Value = 1
$ws.Range("K16").Copy($ws.Range("E30"))
$ws.Range("E30").Value = -100
$ws.Range("C15").Copy($ws.Range("G30"))
$ws.Range("G30").Value = 1
$ws.Range("K16").Copy($ws.Range("H30"))
$ws.Range("H30").Value = -100
$ws.Range("C15").Copy($ws.Range("D31"))
$ws.Range("D31").Value = 2
$ws.Range("K16").Copy($ws.Range("E31"))
$ws.Range("E31").Value = -100
$ws.Range("C15").Copy($ws.Range("G31"))
$ws.Range("G31").Value = 2
$ws.Range("K16").Copy($ws.Range("H31"))
$ws.Range("H31").Value = -100

# --- Plain numeric value updates ---
$ws.Range("F15").Value = 3
$ws.Range("I15").Value = 12
$ws.Range("K15").Value = 100
$ws.Range("L15").Value = 140
$ws.Range("N15").Value = -20
$ws.Range("F16").Value = 2
$ws.Range("H16").Value = -60
$ws.Range("L16").Value = -8.108108108108
$ws.Range("M16").Value = -45.16129032258
$ws.Range("N16").Value = -79.518072289156
$ws.Range("C17").Value = 1
$ws.Range("D17").Value = 4
$ws.Range("E17").Value = -75
$ws.Range("F17").Value = 10
$ws.Range("G17").Value = 12
$ws.Range("H17").Value = -16.666666666666
$ws.Range("I17").Value = 99
$ws.Range("J17").Value = 111
$ws.Range("K17").Value = -10.81081081081
$ws.Range("L17").Value = 57.142857142857
$ws.Range("M17").Value = 20.731707317073
$ws.Range("N17").Value = -43.103448275862
$ws.Range("C18").Value = 2
$ws.Range("D18").Value = 3
$ws.Range("E18").Value = -33.333333333333
$ws.Range("F18").Value = 4
$ws.Range("G18").Value = 11
$ws.Range("H18").Value = -63.636363636363
$ws.Range("I18").Value = 49
$ws.Range("J18").Value = 62
$ws.Range("K18").Value = -20.967741935483
$ws.Range("L18").Value = 25.641025641025
$ws.Range("M18").Value = -57.758620689655
$ws.Range("N18").Value = -93.890274314214
$ws.Range("C19").Value = 7
$ws.Range("D19").Value = 6
$ws.Range("E19").Value = 16.666666666666
$ws.Range("G19").Value = 41
$ws.Range("H19").Value = -24.390243902439
$ws.Range("I19").Value = 297
$ws.Range("J19").Value = 284
$ws.Range("K19").Value = 4.577464788732
$ws.Range("L19").Value = 50
$ws.Range("M19").Value = 23.236514522821
$ws.Range("N19").Value = -39.634146341463
$ws.Range("D20").Value = 4
$ws.Range("E20").Value = -75
$ws.Range("F20").Value = 11
$ws.Range("G20").Value = 18
$ws.Range("H20").Value = -38.888888888888
$ws.Range("I20").Value = 48
$ws.Range("J20").Value = 68
$ws.Range("K20").Value = -29.411764705882
$ws.Range("L20").Value = -36
$ws.Range("M20").Value = -26.153846153846
$ws.Range("N20").Value = -97.085610200364
$ws.Range("C21").Value = 12
$ws.Range("D21").Value = 17
$ws.Range("E21").Value = -29.411764705882
$ws.Range("F21").Value = 61
$ws.Range("G21").Value = 87
$ws.Range("H21").Value = -29.885057471264
$ws.Range("I21").Value = 539
$ws.Range("J21").Value = 564
$ws.Range("K21").Value = -4.432624113475
$ws.Range("L21").Value = 28.947368421052
$ws.Range("M21").Value = -6.74740484429
$ws.Range("N21").Value = -83.661715671415
$ws.Range("F23").Value = 3
$ws.Range("G23").Value = 1
$ws.Range("H23").Value = 200
$ws.Range("I23").Value = 13
$ws.Range("K23").Value = -40.90909090909
$ws.Range("L23").Value = 8.333333333333
$ws.Range("M23").Value = 8.333333333333
$ws.Range("C24").Value = 23
$ws.Range("D24").Value = 29
$ws.Range("E24").Value = -20.689655172413
$ws.Range("F24").Value = 82
$ws.Range("G24").Value = 91
$ws.Range("H24").Value = -9.890109890109
$ws.Range("I24").Value = 620
$ws.Range("J24").Value = 669
$ws.Range("K24").Value = -7.324364723467
$ws.Range("L24").Value = 38.392857142857
$ws.Range("M24").Value = -40.15444015444
$ws.Range("C25").Value = 8
$ws.Range("D25").Value = 12
$ws.Range("E25").Value = -33.333333333333
$ws.Range("F25").Value = 26
$ws.Range("G25").Value = 36
$ws.Range("H25").Value = -27.777777777777
$ws.Range("I25").Value = 283
$ws.Range("J25").Value = 285
$ws.Range("K25").Value = -0.701754385964
$ws.Range("L25").Value = 172.115384615385
$ws.Range("C26").Value = 7
$ws.Range("D26").Value = 8
$ws.Range("E26").Value = -12.5
$ws.Range("F26").Value = 19
$ws.Range("G26").Value = 37
$ws.Range("H26").Value = -48.648648648648
$ws.Range("I26").Value = 214
$ws.Range("J26").Value = 192
$ws.Range("K26").Value = 11.458333333333
$ws.Range("L26").Value = 3.381642512077
$ws.Range("M26").Value = -36.686390532544
$ws.Range("F27").Value = 3
$ws.Range("I27").Value = 15
$ws.Range("K27").Value = 66.666666666666
$ws.Range("L27").Value = 25
$ws.Range("C28").Value = 1
$ws.Range("F28").Value = 5
$ws.Range("H28").Value = 66.666666666666
$ws.Range("I28").Value = 25
$ws.Range("J28").Value = 23
$ws.Range("K28").Value = 8.695652173913
$ws.Range("L28").Value = 19.047619047619
$ws.Range("J29").Value = 2
$ws.Range("K29").Value = -50
$ws.Range("J30").Value = 2
$ws.Range("K30").Value = -50
$ws.Range("J31").Value = 6
$ws.Range("K31").Value = -33.333333333333

Write-Host "applied"
